$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

$ws1.Range("H98").Value = 111112900
$ws1.Range("I98").Value = 125001760
$ws1.Range("K98").Value = 125001760
$ws1.Range("M98").Value = -125000262
$ws1.Range("H116").Value = 6636.1816
$ws1.Range("I116").Value = 6374.75
$ws1.Range("J116").Value = 7333.3335
$ws1.Range("K116").Value = 6374.75
$ws1.Range("L116").Value = 7333.3335
$ws1.Range("M116").Value = -2932.75
$ws1.Range("N116").Value = -14217.3335
$ws1.Range("H122").Value = 111112900
$ws1.Range("I122").Value = 125001760
$ws1.Range("K122").Value = 375005280
$ws1.Range("M122").Value = -375002830
$ws1.Range("H137").Value = 3529.75
$ws1.Range("J137").Value = 3636.35
$ws1.Range("L137").Value = 10909.05
$ws1.Range("N137").Value = -16009.05
$ws1.Range("H138").Value = 2349.1506
$ws1.Range("I138").Value = 1285.6072
$ws1.Range("K138").Value = 3856.8216
$ws1.Range("M138").Value = 1283.1784
$ws2.Range("H32").Value = 7249242.5
$ws2.Range("I32").Value = 7938872.5
$ws2.Range("K32").Value = 7938872.5
$ws2.Range("M32").Value = -7938585.5
$ws2.Range("H63").Value = 5037.625
$ws2.Range("I63").Value = 5413.5713
$ws2.Range("K63").Value = 5413.5713
$ws2.Range("M63").Value = -4727.5713
$ws2.Range("H66").Value = 5037.625
$ws2.Range("I66").Value = 5413.5713
$ws2.Range("K66").Value = 27067.8565
$ws2.Range("M66").Value = -23635.8565
$ws2.Range("H74").Value = 7583713.5
$ws2.Range("I74").Value = 10418868
$ws2.Range("K74").Value = 10418868
$ws2.Range("M74").Value = -10417994
$ws2.Range("H77").Value = 7583713.5
$ws2.Range("I77").Value = 10418868
$ws2.Range("K77").Value = 52094340
$ws2.Range("M77").Value = -52089972
$ws2.Range("H131").Value = 65238.332
$ws2.Range("J131").Value = 65238.332
$ws2.Range("L131").Value = 65238.332
$ws2.Range("N131").Value = -75318.33199999999
$ws2.Range("H132").Value = 2811.093
$ws2.Range("I132").Value = 1059.0667
$ws2.Range("J132").Value = 6854.231
$ws2.Range("K132").Value = 3177.2001
$ws2.Range("L132").Value = 20562.693
$ws2.Range("M132").Value = -647.2001
$ws2.Range("N132").Value = -25622.693
$ws3.Range("H63").Value = 118499.5
$ws3.Range("J63").Value = 118499.5
$ws3.Range("L63").Value = 118499.5
$ws3.Range("N63").Value = -119871.5
$ws3.Range("H66").Value = 118499.5
$ws3.Range("J66").Value = 118499.5
$ws3.Range("L66").Value = 355498.5
$ws3.Range("N66").Value = -362362.5
$ws3.Range("H99").Value = 5646.355
$ws3.Range("I99").Value = 6544.5264
$ws3.Range("K99").Value = 6544.5264
$ws3.Range("M99").Value = -5046.5264
$ws3.Range("H105").Value = 2483.9
$ws3.Range("I105").Value = 1961.6
$ws3.Range("K105").Value = 1961.6
$ws3.Range("M105").Value = -214.5999999999999
$ws3.Range("H107").Value = 1841.2258
$ws3.Range("I107").Value = 1473.5416
$ws3.Range("K107").Value = 1473.5416
$ws3.Range("M107").Value = 446.4584
$ws3.Range("H133").Value = 55999.97
$ws3.Range("I133").Value = 19000
$ws3.Range("J133").Value = 57121.184
$ws3.Range("K133").Value = 19000
$ws3.Range("L133").Value = 57121.184
$ws3.Range("M133").Value = -13940
$ws3.Range("N133").Value = -67241.18400000001
$ws3.Range("H135").Value = 57878.758
$ws3.Range("J135").Value = 57878.758
$ws3.Range("L135").Value = 57878.758
$ws3.Range("N135").Value = -68018.758
$ws4.Range("H31").Value = 525999.5
$ws4.Range("I31").Value = 10031.546
$ws4.Range("J31").Value = 1019534.06
$ws4.Range("K31").Value = 10031.546
$ws4.Range("L31").Value = 1019534.06
$ws4.Range("M31").Value = -9736.546
$ws4.Range("N31").Value = -1020124.06
$ws4.Range("H34").Value = 525999.5
$ws4.Range("I34").Value = 10031.546
$ws4.Range("J34").Value = 1019534.06
$ws4.Range("K34").Value = 10031.546
$ws4.Range("L34").Value = 1019534.06
$ws4.Range("M34").Value = -9829.546
$ws4.Range("N34").Value = -1019938.06
$ws4.Range("H50").Value = 39999
$ws4.Range("J50").Value = 0
$ws4.Range("L50").Value = 0
$ws4.Range("N50").Value = $null
$ws4.Range("H86").Value = 5828
$ws4.Range("J86").Value = 5449
$ws4.Range("L86").Value = 5449
$ws4.Range("N86").Value = -7695
$ws4.Range("H89").Value = 5828
$ws4.Range("J89").Value = 5449
$ws4.Range("L89").Value = 27245
$ws4.Range("N89").Value = -38477
$ws4.Range("H132").Value = 2786.7188
$ws4.Range("I132").Value = 2006.1538
$ws4.Range("K132").Value = 6018.4614
$ws4.Range("M132").Value = -3488.4614
$ws5.Range("H41").Value = 2757.7144
$ws5.Range("I41").Value = 55
$ws5.Range("J41").Value = 3208.1667
$ws5.Range("K41").Value = 165
$ws5.Range("L41").Value = 9624.500100000001
$ws5.Range("M41").Value = 173
$ws5.Range("N41").Value = -10300.5001
$ws5.Range("H131").Value = 8250.462
$ws5.Range("J131").Value = 7373.3184
$ws5.Range("L131").Value = 22119.9552
$ws5.Range("N131").Value = -32199.9552
$ws6.Range("H80").Value = 3067.4614
$ws6.Range("I80").Value = 2874.6667
$ws6.Range("J80").Value = 3501.25
$ws6.Range("K80").Value = 2874.6667
$ws6.Range("L80").Value = 3501.25
$ws6.Range("M80").Value = -1876.6667
$ws6.Range("N80").Value = -5497.25
$ws6.Range("H83").Value = 3067.4614
$ws6.Range("I83").Value = 2874.6667
$ws6.Range("J83").Value = 3501.25
$ws6.Range("K83").Value = 14373.3335
$ws6.Range("L83").Value = 17506.25
$ws6.Range("M83").Value = -9381.333500000001
$ws6.Range("N83").Value = -27490.25
$ws6.Range("H102").Value = 3199.7097
$ws6.Range("I102").Value = 2473.0334
$ws6.Range("K102").Value = 2473.0334
$ws6.Range("M102").Value = -851.0333999999998
$ws6.Range("H132").Value = 142873100
$ws6.Range("I132").Value = 166668600
$ws6.Range("K132").Value = 500005800
$ws6.Range("M132").Value = -500003270
$ws7.Range("H22").Value = 2122
$ws7.Range("I22").Value = 1490
$ws7.Range("J22").Value = 2332.6667
$ws7.Range("K22").Value = 1490
$ws7.Range("L22").Value = 2332.6667
$ws7.Range("M22").Value = -1195
$ws7.Range("N22").Value = -2922.6667
$ws7.Range("H27").Value = 2122
$ws7.Range("I27").Value = 1490
$ws7.Range("J27").Value = 2332.6667
$ws7.Range("K27").Value = 1490
$ws7.Range("L27").Value = 2332.6667
$ws7.Range("M27").Value = -1383
$ws7.Range("N27").Value = -2546.6667
$ws7.Range("H40").Value = 4665
$ws7.Range("I40").Value = 3663.875
$ws7.Range("K40").Value = 3663.875
$ws7.Range("M40").Value = -3527.875
$ws7.Range("H46").Value = 4927.5713
$ws7.Range("J46").Value = 7354.1816
$ws7.Range("L46").Value = 7354.1816
$ws7.Range("N46").Value = -7730.1816
$ws7.Range("H55").Value = 45454972
$ws7.Range("J55").Value = 444.14285
$ws7.Range("L55").Value = 444.14285
$ws7.Range("N55").Value = -790.14285
$ws7.Range("H58").Value = 18194.4
$ws7.Range("I58").Value = 13500
$ws7.Range("J58").Value = 21324
$ws7.Range("K58").Value = 13500
$ws7.Range("L58").Value = 21324
$ws7.Range("M58").Value = -13240
$ws7.Range("N58").Value = -21844
$ws7.Range("H82").Value = 737.7273
$ws7.Range("I82").Value = 611.5
$ws7.Range("J82").Value = 2000
$ws7.Range("K82").Value = 611.5
$ws7.Range("L82").Value = 2000
$ws7.Range("M82").Value = -250.5
$ws7.Range("N82").Value = -2722
$ws7.Range("H85").Value = 737.7273
$ws7.Range("I85").Value = 611.5
$ws7.Range("J85").Value = 2000
$ws7.Range("K85").Value = 611.5
$ws7.Range("L85").Value = 2000
$ws7.Range("M85").Value = 636.5
$ws7.Range("N85").Value = -4496
$ws7.Range("H93").Value = 55563676
$ws7.Range("I93").Value = 66675956
$ws7.Range("K93").Value = 66675956
$ws7.Range("M93").Value = -66674708
$ws7.Range("H122").Value = 5811.735
$ws7.Range("I122").Value = 4250.5815
$ws7.Range("J122").Value = 17000
$ws7.Range("K122").Value = 12751.7445
$ws7.Range("L122").Value = 51000
$ws7.Range("M122").Value = -10301.7445
$ws7.Range("N122").Value = -55900
$ws8.Range("H81").Value = 13000
$ws8.Range("I81").Value = 5000
$ws8.Range("K81").Value = 10000
$ws8.Range("M81").Value = -8939
$ws8.Range("H84").Value = 13000
$ws8.Range("I84").Value = 5000
$ws8.Range("K84").Value = 50000
$ws8.Range("M84").Value = -44696
$ws8.Range("H100").Value = 1815
$ws8.Range("I100").Value = 1798.1818
$ws8.Range("K100").Value = 3596.3636
$ws8.Range("M100").Value = -3055.3636
$ws8.Range("H132").Value = 2413.05
$ws8.Range("I132").Value = 2297.8125
$ws8.Range("K132").Value = 6893.4375
$ws8.Range("M132").Value = -4363.4375
